$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6332324147224426
$ws.Range("B1").Value = 1.977708220481873
$ws.Range("C1").Value = 3.415624380111694
$ws.Range("D1").Value = 1.724713325500488
$ws.Range("E1").Value = 0.748962938785553
